$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (changed) date column for rows 2-6 from 45184 (2023-09-15)
# to 45185 (2023-09-16), preserving existing formatting.
foreach ($r in 2..6) {
    $ws.Cells.Item($r, 3).Value = 45185
}
